# Registration TestCase workbook update
# - Adds a new "TC009 - Verify duplicate registration with same face" test case
#   (rows 17-18) to the existing test-case table.
# - Widens the "Expected Result" column area (D:E) to better fit the new text.
# - Normalizes row heights (auto-fit) for the rows that are re-laid out as a
#   result of the edit, including collapsing the blank spacer row (row 9)
#   back to an implicit (un-stored) row.
# - Leaves the final selection on D18 (the last cell typed), matching the
#   cursor position the author ended up at.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Add the new test case content (TC009) in rows 17-18.
#    Column order below intentionally enters Steps (column C) for both
#    new rows before Expected Result (column D) for both new rows, which
#    is the order the content was actually typed in.
# ---------------------------------------------------------------------
$ws.Range("A17").Value2 = "TC009"
$ws.Range("B17").Value2 = "Verify duplicate registration with same face"
$ws.Range("C17").Value2 = "1. Enter registration fields"
$ws.Range("C18").Value2 = "2. Click Add user"
$ws.Range("D17").Value2 = "The system should not allow registring same face.(Pass)"
$ws.Range("D18").Value2 = "This can allo same face if camera quality is low, lightening issue(Fail)"

# ---------------------------------------------------------------------
# 2. Widen column D/E so the longer Expected Result text fits.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 55.83
$ws.Columns.Item(5).ColumnWidth = 55.83

# ---------------------------------------------------------------------
# 3. Re-layout (auto fit) the rows whose explicit row height is no longer
#    needed, and collapse the blank row 9 back down to the sheet default.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).AutoFit()
$ws.Rows.Item(16).AutoFit()
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).AutoFit()

# ---------------------------------------------------------------------
# 4. Update the view: scroll so column B is left-most visible, and leave
#    the active selection on the last cell that was edited (D18).
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("D18").Select()
